$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (I1, J1) - copy style from existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data rows 2-15: I column = 1, J column = same value as H column
for ($r = 2; $r -le 15; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
